$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: set to former row 9 data
$ws.Range("D2").Value = 44489
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 24000
$ws.Range("O2").Value = 25000
$ws.Range("P2").Value = 24500
$ws.Range("Q2").Value = '$/caja 12 kilos'
$ws.Range("S2").Value = 2042
$ws.Range("T2").Value = 12

# Row 3: set to former row 2 data
$ws.Range("D3").Value = 44468
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 29000
$ws.Range("O3").Value = 30000
$ws.Range("P3").Value = 29500
$ws.Range("Q3").Value = '$/bandeja 10 kilos'
$ws.Range("S3").Value = 2950
$ws.Range("T3").Value = 10

# Row 5: set to former row 11 data
$ws.Range("D5").Value = 44160
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 19000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 19500
$ws.Range("Q5").Value = '$/caja 13 kilos'
$ws.Range("S5").Value = 1500
$ws.Range("T5").Value = 13

# Row 6: set to former row 8 data
$ws.Range("D6").Value = 44524
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 23000
$ws.Range("O6").Value = 24000
$ws.Range("P6").Value = 23500
$ws.Range("Q6").Value = '$/caja 12 kilos'
$ws.Range("S6").Value = 1958
$ws.Range("T6").Value = 12

# Row 7: set to former row 5 data
$ws.Range("D7").Value = 44441
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 29000
$ws.Range("O7").Value = 30000
$ws.Range("P7").Value = 29500
$ws.Range("Q7").Value = '$/caja 12 kilos'
$ws.Range("S7").Value = 2458
$ws.Range("T7").Value = 12

# Row 8: set to former row 10 data
$ws.Range("D8").Value = 44167
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 18000
$ws.Range("O8").Value = 19000
$ws.Range("P8").Value = 18500
$ws.Range("Q8").Value = '$/caja 13 kilos'
$ws.Range("S8").Value = 1423
$ws.Range("T8").Value = 13

# Row 9: set to former row 6 data
$ws.Range("D9").Value = 44475
$ws.Range("L9").Value = "Especial"
$ws.Range("M9").Value = 200
$ws.Range("N9").Value = 32000
$ws.Range("O9").Value = 33000
$ws.Range("P9").Value = 32500
$ws.Range("Q9").Value = '$/caja 12 kilos'
$ws.Range("S9").Value = 2708
$ws.Range("T9").Value = 12

# Row 10: set to former row 3 data
$ws.Range("D10").Value = 44496
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 23000
$ws.Range("O10").Value = 24000
$ws.Range("P10").Value = 23500
$ws.Range("Q10").Value = '$/caja 12 kilos'
$ws.Range("S10").Value = 1958
$ws.Range("T10").Value = 12

# Row 11: set to former row 7 data
$ws.Range("D11").Value = 44482
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 160
$ws.Range("N11").Value = 25000
$ws.Range("O11").Value = 26000
$ws.Range("P11").Value = 25500
$ws.Range("Q11").Value = '$/caja 12 kilos'
$ws.Range("S11").Value = 2125
$ws.Range("T11").Value = 12
